$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.347.60"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.76%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.427.63"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.45%  "

# Row 4
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$ws.Range("B5").Value = "Solana"
$ws.Range("C5").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "156.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.37%  "

# Row 6
$ws.Range("B6").Value = "BNB"
$ws.Range("C6").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "490.82"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.86%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.995"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.54%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.606"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +20.40%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.440.77"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.26%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.27"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +11.97%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.101"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.84%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.333"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.25%  "

# Row 13
$ws.Range("E13").Value = "  +1.47%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.832.68"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.21%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "57.327.06"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.34%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.86"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.16%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000134"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.11%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.446.82"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.26%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.69"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.14%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "323.29"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.89%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.07"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.47%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.996"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.10%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.92"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.19%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "58.31"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.32%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.404"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.17%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.991"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.07%  "

# Row 27
$ws.Range("E27").Value = "  -0.43%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.522.73"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.00%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.33"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.42%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0801"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.42%  "

# Row 31
$ws.Range("E31").Value = "  -0.19%  "

# Row 32
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.81"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.14%  "

# Row 33
$ws.Range("B33").Value = "Monero"
$ws.Range("C33").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "150.42"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.80%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.53"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.75%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.35"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.53%  "

# Row 36
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.15"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.90%  "

# Row 37
$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.77"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.56%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.824"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.88%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "34.39"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.66%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.39"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.89%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.55"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.55%  "

# Row 42
$ws.Range("E42").Value = "  +5.70%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "278.84"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.20%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.992"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.60%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.596"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.09%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0538"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.12%  "

# Row 47
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0230"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.56%  "

# Row 48
$ws.Range("B48").Value = "WhiteBITCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.24"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.10%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.66"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.05%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.98"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.66%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.691"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +9.86%  "

